$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.589944
$ws.Range("H2").Value = 1.769832
$ws.Range("I2").Value = 0.1315369184485108
$ws.Range("J2").Value = 0.1315369184485108
$ws.Range("Q2").Value = 1.121454440096
$ws.Range("R2").Value = 10.093089960864
$ws.Range("S2").Value = 0.1315369184485108
$ws.Range("T2").Value = 0.1315369184485108

# Row 3
$ws.Range("I3").Value = 0.1512402693388587
$ws.Range("J3").Value = 0.1512402693388587
$ws.Range("S3").Value = 0.1512402693388587
$ws.Range("T3").Value = 0.1512402693388587

# Row 4
$ws.Range("G4").Value = 0.3049506666666666
$ws.Range("H4").Value = 0.914852
$ws.Range("I4").Value = 0.06799335355924008
$ws.Range("J4").Value = 0.06799335355924008
$ws.Range("Q4").Value = 0.5796961731004444
$ws.Range("R4").Value = 5.217265557904
$ws.Range("S4").Value = 0.06799335355924008
$ws.Range("T4").Value = 0.06799335355924008

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4130943333333333
$ws.Range("H5").Value = 1.239283
$ws.Range("I5").Value = 0.09210561618595764
$ws.Range("J5").Value = 0.09210561618595764
$ws.Range("Q5").Value = 0.7852719483462222
$ws.Range("R5").Value = 7.067447535115999
$ws.Range("S5").Value = 0.09210561618595764
$ws.Range("T5").Value = 0.09210561618595764

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.8655823333333333
$ws.Range("H6").Value = 2.596747
$ws.Range("I6").Value = 0.1929946448987334
$ws.Range("J6").Value = 0.1929946448987334
$ws.Range("Q6").Value = 1.645429313604889
$ws.Range("R6").Value = 14.808863822444
$ws.Range("S6").Value = 0.1929946448987334
$ws.Range("T6").Value = 0.1929946448987334

# Row 7
$ws.Range("G7").Value = 1.633122
$ws.Range("H7").Value = 4.899366000000001
$ws.Range("I7").Value = 0.3641291975686995
$ws.Range("J7").Value = 0.3641291975686995
$ws.Range("Q7").Value = 3.104484354648001
$ws.Range("R7").Value = 27.940359191832
$ws.Range("S7").Value = 0.3641291975686995
$ws.Range("T7").Value = 0.3641291975686995
